$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "25.811.29"
Set-TextValue "D3" "1.816.22"
Set-TextValue "E3" "  -3.07%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "276.86"
Set-TextValue "E5" "  -8.02%  "
Set-TextValue "E6" "  -0.04%  "
Set-TextValue "D7" "0.5096"
Set-TextValue "E7" "  -4.55%  "
Set-TextValue "D8" "0.3532"
Set-TextValue "E8" "  -6.09%  "
Set-TextValue "D9" "44.65"
Set-TextValue "E9" "  -1.86%  "
Set-TextValue "D10" "0.06667"
Set-TextValue "E10" "  -7.17%  "
Set-TextValue "D11" "20.04"
Set-TextValue "E11" "  -7.32%  "
Set-TextValue "D12" "0.8292"
Set-TextValue "E12" "  -6.54%  "
Set-TextValue "D13" "0.07861"
Set-TextValue "E13" "  -3.79%  "
Set-TextValue "D14" "1.812.65"
Set-TextValue "E14" "  -3.11%  "
Set-TextValue "D15" "5.081"
Set-TextValue "E15" "  -3.81%  "
Set-TextValue "D16" "87.54"
Set-TextValue "E16" "  -6.39%  "
Set-TextValue "D17" "1.000"
Set-TextValue "E17" "  +0.03%  "
Set-TextValue "E18" "  -4.29%  "
Set-TextValue "D19" "0.000008052"
Set-TextValue "E19" "  -5.83%  "
Set-TextValue "D20" "1.000"
Set-TextValue "E20" "  -0.04%  "
Set-TextValue "D21" "25.852.76"
Set-TextValue "E21" "  -3.69%  "
Set-TextValue "D22" "4.732"
Set-TextValue "E22" "  -5.05%  "
Set-TextValue "D23" "10.01"
Set-TextValue "E23" "  -6.24%  "
Set-TextValue "D24" "6.092"
Set-TextValue "E24" "  -4.70%  "
Set-TextValue "B25" "LidoDAOToken"
Set-TextValue "C25" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D25" "2.201"
Set-TextValue "E25" "  -2.96%  "
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "141.19"
Set-TextValue "E26" "  -3.59%  "
Set-TextValue "E27" "  -3.63%  "
Set-TextValue "D28" "17.09"
Set-TextValue "E28" "  -5.20%  "
Set-TextValue "D29" "109.47"
Set-TextValue "E29" "  -3.81%  "
Set-TextValue "E30" "  -7.87%  "
Set-TextValue "D31" "4.245"
Set-TextValue "E31" "  -7.97%  "
Set-TextValue "D32" "0.08799"
Set-TextValue "E32" "  -3.78%  "
Set-TextValue "D33" "0.04905"
Set-TextValue "E33" "  -1.38%  "
Set-TextValue "D34" "0.7289"
Set-TextValue "E35" "  -3.07%  "
Set-TextValue "D36" "2.864"
Set-TextValue "E36" "  -3.41%  "
Set-TextValue "D37" "0.9999"
Set-TextValue "E37" "  -0.12%  "
Set-TextValue "E38" "  -1.93%  "
Set-TextValue "D39" "2.375"
Set-TextValue "E39" "  -8.23%  "
Set-TextValue "B40" "TheSandbox"
Set-TextValue "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D40" "0.5197"
Set-TextValue "E40" "  -13.66%  "
Set-TextValue "B41" "VeChain"
Set-TextValue "C41" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.01852"
Set-TextValue "E41" "  -4.97%  "
Set-TextValue "D42" "0.9569"
Set-TextValue "E42" "  -10.74%  "
Set-TextValue "D43" "6.211"
Set-TextValue "E43" "  -5.41%  "
Set-TextValue "D44" "111.10"
Set-TextValue "E44" "  -3.12%  "
Set-TextValue "D45" "8.012"
Set-TextValue "E45" "  -9.73%  "
Set-TextValue "E46" "  -0.01%  "
Set-TextValue "D47" "0.4578"
Set-TextValue "E47" "  -10.97%  "
Set-TextValue "D48" "0.1366"
Set-TextValue "E48" "  -8.32%  "
Set-TextValue "D49" "36.61"
Set-TextValue "E49" "  -2.46%  "
Set-TextValue "D50" "9.238"
Set-TextValue "E50" "  -6.67%  "
Set-TextValue "E51" "  -7.94%  "

Write-Host "Applied 97 cell updates"
